$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of row 11 (the last existing data row) down into the
# new row 12, so new cells inherit the right borders/styles.
$ws.Range("A11:F11").Copy()
$ws.Range("A12:F12").PasteSpecial(-4122)

# Populate the new test case row.
$ws.Range("A12").Value = "FLD_Transmittals_ActionRequired_CaC_IssuedForReview"
$ws.Range("B12").Value = "Verifies the Close and Cancel option for the user in the Transmital record"
$ws.Range("C12").Value = "N"
$ws.Range("D12").Value = "Y"
$ws.Range("E12").Value = "PASS"
$ws.Range("F12").Value = "Sprint2"

# B12 gets its own alignment (left/top, no wrap) distinct from column A.
$ws.Range("B12").HorizontalAlignment = -4131
$ws.Range("B12").VerticalAlignment = -4160
$ws.Range("B12").WrapText = $false

# Match the row height used by the other multi-line rows.
$ws.Rows.Item(12).RowHeight = 30

# Extend the existing validation lists to cover the new row.
$ws.Range("C2:D11").Validation.Delete()
$ws.Range("C2:D12").Validation.Add(3, 1, 1, """Y,N""")
$ws.Range("F2:F11").Validation.Delete()
$ws.Range("F2:F12").Validation.Add(3, 1, 1, """Sprint1,Sprint2,Sprint3,Sprint4,Sprint5,Sprint6,Sprint7,Sprint8,Sprint9,Sprint10""")

# Select the new row's first cell, matching the saved cursor position.
[void]$ws.Range("A12").Select()

Write-Host "done"
